$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 157; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2
    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "
        if ($parts.Count -gt 1) {
            $reversed = $parts[($parts.Count - 1)..0]
            $newVal = [string]::Join(", ", $reversed)
            $cell.Value2 = $newVal
        }
    }
}
